$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 493.16666
$ws.Range("J32").Value = 512.6429000000001
$ws.Range("L32").Value = 512.6429000000001
$ws.Range("N32").Value = -1164.6429
$ws.Range("H132").Value = 4764353
$ws.Range("I132").Value = 5104449.5
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 15313348.5
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -15310818.5
$ws.Range("N132").Value = -14058.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21123.37
$ws.Range("I32").Value = 23232.36
$ws.Range("J32").Value = 6059.143
$ws.Range("K32").Value = 23232.36
$ws.Range("L32").Value = 6059.143
$ws.Range("M32").Value = -22945.36
$ws.Range("N32").Value = -6633.143
$ws.Range("H61").Value = 2369.2307
$ws.Range("I61").Value = 971.4286
$ws.Range("K61").Value = 971.4286
$ws.Range("M61").Value = -759.4286
$ws.Range("H97").Value = 715.88464
$ws.Range("I97").Value = 680.6842
$ws.Range("J97").Value = 811.4286
$ws.Range("K97").Value = 680.6842
$ws.Range("L97").Value = 811.4286
$ws.Range("M97").Value = -184.6842
$ws.Range("N97").Value = -1803.4286
$ws.Range("H136").Value = 2369.2307
$ws.Range("I136").Value = 971.4286
$ws.Range("K136").Value = 2914.2858
$ws.Range("M136").Value = -364.2857999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 21542.666
$ws.Range("J76").Value = 21542.666
$ws.Range("L76").Value = 21542.666
$ws.Range("N76").Value = -22172.666
$ws.Range("H79").Value = 21542.666
$ws.Range("J79").Value = 21542.666
$ws.Range("L79").Value = 21542.666
$ws.Range("N79").Value = -23726.666
$ws.Range("H134").Value = 38064.605
$ws.Range("I134").Value = 52115.5
$ws.Range("J134").Value = 2937.375
$ws.Range("K134").Value = 156346.5
$ws.Range("L134").Value = 8812.125
$ws.Range("M134").Value = -153811.5
$ws.Range("N134").Value = -13882.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 104.166664
$ws.Range("I7").Value = 75
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 75
$ws.Range("L7").Value = 250
$ws.Range("M7").Value = 38
$ws.Range("N7").Value = -476
$ws.Range("H31").Value = 8337573.5
$ws.Range("I31").Value = 2941.0952
$ws.Range("J31").Value = 66680000
$ws.Range("K31").Value = 2941.0952
$ws.Range("L31").Value = 66680000
$ws.Range("M31").Value = -2646.0952
$ws.Range("N31").Value = -66680590
$ws.Range("H34").Value = 8337573.5
$ws.Range("I34").Value = 2941.0952
$ws.Range("J34").Value = 66680000
$ws.Range("K34").Value = 2941.0952
$ws.Range("L34").Value = 66680000
$ws.Range("M34").Value = -2739.0952
$ws.Range("N34").Value = -66680404
$ws.Range("H99").Value = 1259.1364
$ws.Range("I99").Value = 895.1
$ws.Range("J99").Value = 1562.5
$ws.Range("K99").Value = 895.1
$ws.Range("L99").Value = 1562.5
$ws.Range("M99").Value = 602.9
$ws.Range("N99").Value = -4558.5
$ws.Range("H126").Value = 1259.1364
$ws.Range("I126").Value = 895.1
$ws.Range("J126").Value = 1562.5
$ws.Range("K126").Value = 2685.3
$ws.Range("L126").Value = 4687.5
$ws.Range("M126").Value = -215.3000000000002
$ws.Range("N126").Value = -9627.5
$ws.Range("H132").Value = 1926.9333
$ws.Range("I132").Value = 1304.5834
$ws.Range("J132").Value = 4416.3335
$ws.Range("K132").Value = 3913.7502
$ws.Range("L132").Value = 13249.0005
$ws.Range("M132").Value = -1383.7502
$ws.Range("N132").Value = -18309.0005
$ws.Range("H134").Value = 1025.3077
$ws.Range("I134").Value = 889.4783
$ws.Range("K134").Value = 2668.4349
$ws.Range("M134").Value = -133.4349000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 13895462
$ws.Range("J9").Value = 14499438
$ws.Range("L9").Value = 43498314
$ws.Range("N9").Value = -43498762
$ws.Range("H41").Value = 975
$ws.Range("I41").Value = 950
$ws.Range("J41").Value = 1000
$ws.Range("K41").Value = 2850
$ws.Range("L41").Value = 3000
$ws.Range("M41").Value = -2512
$ws.Range("N41").Value = -3676
$ws.Range("H69").Value = 1495.9048
$ws.Range("J69").Value = 1495.9048
$ws.Range("L69").Value = 4487.7144
$ws.Range("N69").Value = -6109.7144
$ws.Range("H72").Value = 1495.9048
$ws.Range("J72").Value = 1495.9048
$ws.Range("L72").Value = 13463.1432
$ws.Range("N72").Value = -21575.1432
$ws.Range("H93").Value = 4818.1816
$ws.Range("I93").Value = 1000
$ws.Range("K93").Value = 3000
$ws.Range("M93").Value = -1128
$ws.Range("H132").Value = 1814.2858
$ws.Range("I132").Value = 500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 44.4
$ws.Range("I2").Value = 21.333334
$ws.Range("J2").Value = 59.77778
$ws.Range("K2").Value = 21.333334
$ws.Range("L2").Value = 59.77778
$ws.Range("M2").Value = 91.66666599999999
$ws.Range("N2").Value = -285.77778
$ws.Range("H102").Value = 1447.875
$ws.Range("I102").Value = 1437.8462
$ws.Range("K102").Value = 1437.8462
$ws.Range("M102").Value = 184.1538
$ws.Range("H126").Value = 8785.714
$ws.Range("I126").Value = 7000
$ws.Range("J126").Value = 13250
$ws.Range("K126").Value = 21000
$ws.Range("L126").Value = 39750
$ws.Range("M126").Value = -18530
$ws.Range("N126").Value = -44690
$ws.Range("H132").Value = 93501.63
$ws.Range("I132").Value = 113124.336
$ws.Range("K132").Value = 339373.008
$ws.Range("M132").Value = -336843.008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3003998.8
$ws.Range("I93").Value = 3003998.8
$ws.Range("K93").Value = 3003998.8
$ws.Range("M93").Value = -3002750.8
$ws.Range("H100").Value = 2798.8
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459
$ws.Range("H122").Value = 2347.8096
$ws.Range("I122").Value = 2086.9333
$ws.Range("K122").Value = 6260.7999
$ws.Range("M122").Value = -3810.7999
$ws.Range("H132").Value = 10483.458
$ws.Range("I132").Value = 13394.294
$ws.Range("J132").Value = 3414.2856
$ws.Range("K132").Value = 40182.882
$ws.Range("L132").Value = 10242.8568
$ws.Range("M132").Value = -37652.882
$ws.Range("N132").Value = -15302.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 52194
$ws.Range("J123").Value = 52194
$ws.Range("L123").Value = 52194
$ws.Range("N123").Value = -61994
$ws.Range("H136").Value = 6882.2856
$ws.Range("I136").Value = 8048.706
$ws.Range("K136").Value = 24146.118
$ws.Range("M136").Value = -21596.118
